$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 3039
$ws.Range("I20").Value = 3039
$ws.Range("K20").Value = 3039
$ws.Range("M20").Value = -2809
$ws.Range("H35").Value = 3039
$ws.Range("I35").Value = 3039
$ws.Range("K35").Value = 3039
$ws.Range("M35").Value = -2660
$ws.Range("H40").Value = 4031.6667
$ws.Range("I40").Value = 2101
$ws.Range("K40").Value = 2101
$ws.Range("M40").Value = -1926
$ws.Range("H74").Value = 4355.3335
$ws.Range("I74").Value = 4355.3335
$ws.Range("K74").Value = 4355.3335
$ws.Range("M74").Value = -3419.3335
$ws.Range("H77").Value = 4355.3335
$ws.Range("I77").Value = 4355.3335
$ws.Range("K77").Value = 21776.6675
$ws.Range("M77").Value = -17096.6675
$ws.Range("H86").Value = 5200.6
$ws.Range("J86").Value = 4667.6665
$ws.Range("L86").Value = 4667.6665
$ws.Range("N86").Value = -6913.6665
$ws.Range("H89").Value = 5200.6
$ws.Range("J89").Value = 4667.6665
$ws.Range("L89").Value = 23338.3325
$ws.Range("N89").Value = -34570.3325
$ws.Range("H100").Value = 3107.6
$ws.Range("I100").Value = 2879.6667
$ws.Range("J100").Value = 3449.5
$ws.Range("K100").Value = 2879.6667
$ws.Range("L100").Value = 3449.5
$ws.Range("M100").Value = -2338.6667
$ws.Range("N100").Value = -4531.5
$ws.Range("H116").Value = 5749.5
$ws.Range("I116").Value = 4999
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 4999
$ws.Range("L116").Value = 6500
$ws.Range("M116").Value = -1557
$ws.Range("N116").Value = -13384
$ws.Range("H132").Value = 5053.143
$ws.Range("I132").Value = 5245.846
$ws.Range("J132").Value = 2548
$ws.Range("K132").Value = 15737.538
$ws.Range("L132").Value = 7644
$ws.Range("M132").Value = -13207.538
$ws.Range("N132").Value = -12704
$ws.Range("H137").Value = 2438.5
$ws.Range("I137").Value = 2438.5
$ws.Range("K137").Value = 7315.5
$ws.Range("M137").Value = -4765.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2932.7036
$ws.Range("I32").Value = 2660.923
$ws.Range("K32").Value = 2660.923
$ws.Range("M32").Value = -2373.923
$ws.Range("H45").Value = 1605.25
$ws.Range("I45").Value = 1526.9
$ws.Range("K45").Value = 1526.9
$ws.Range("M45").Value = -1149.9
$ws.Range("H102").Value = 1865.2222
$ws.Range("I102").Value = 1865.2222
$ws.Range("K102").Value = 1865.2222
$ws.Range("M102").Value = -243.2221999999999
$ws.Range("H122").Value = 1535.0416
$ws.Range("I122").Value = 1588.7391
$ws.Range("K122").Value = 4766.2173
$ws.Range("M122").Value = -2316.2173
$ws.Range("H132").Value = 2400.8
$ws.Range("I132").Value = 2400.8
$ws.Range("K132").Value = 7202.400000000001
$ws.Range("M132").Value = -4672.400000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2303.2942
$ws.Range("I20").Value = 1323.75
$ws.Range("K20").Value = 1323.75
$ws.Range("M20").Value = -1076.75
$ws.Range("H86").Value = 2901.9412
$ws.Range("I86").Value = 2901.9412
$ws.Range("K86").Value = 2901.9412
$ws.Range("M86").Value = -1778.9412
$ws.Range("H89").Value = 2901.9412
$ws.Range("I89").Value = 2901.9412
$ws.Range("K89").Value = 14509.706
$ws.Range("M89").Value = -8893.706000000002
$ws.Range("H105").Value = 3434.0386
$ws.Range("I105").Value = 2691.625
$ws.Range("J105").Value = 4621.9
$ws.Range("K105").Value = 2691.625
$ws.Range("L105").Value = 4621.9
$ws.Range("M105").Value = -944.625
$ws.Range("N105").Value = -8115.9
$ws.Range("H107").Value = 781.5
$ws.Range("I107").Value = 578.5
$ws.Range("K107").Value = 578.5
$ws.Range("M107").Value = 1341.5
$ws.Range("H134").Value = 5150
$ws.Range("I134").Value = 5600
$ws.Range("K134").Value = 16800
$ws.Range("M134").Value = -14265

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1498.9
$ws.Range("I31").Value = 1498.7778
$ws.Range("K31").Value = 1498.7778
$ws.Range("M31").Value = -1203.7778
$ws.Range("H34").Value = 1498.9
$ws.Range("I34").Value = 1498.7778
$ws.Range("K34").Value = 1498.7778
$ws.Range("M34").Value = -1296.7778
$ws.Range("H62").Value = 3447.25
$ws.Range("J62").Value = 3447.25
$ws.Range("L62").Value = 3447.25
$ws.Range("N62").Value = -4695.25
$ws.Range("H65").Value = 3447.25
$ws.Range("J65").Value = 3447.25
$ws.Range("L65").Value = 17236.25
$ws.Range("N65").Value = -23476.25
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H94").Value = 1178.8
$ws.Range("J94").Value = 681.3333
$ws.Range("L94").Value = 681.3333
$ws.Range("N94").Value = -1583.3333
$ws.Range("H99").Value = 2930
$ws.Range("I99").Value = 2910.8572
$ws.Range("K99").Value = 2910.8572
$ws.Range("M99").Value = -1412.8572
$ws.Range("H105").Value = 1622.75
$ws.Range("I105").Value = 1521.4286
$ws.Range("K105").Value = 1521.4286
$ws.Range("M105").Value = 225.5714
$ws.Range("H107").Value = 1181.25
$ws.Range("I107").Value = 420.25
$ws.Range("K107").Value = 420.25
$ws.Range("M107").Value = 1499.75
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -104909
$ws.Range("H126").Value = 2930
$ws.Range("I126").Value = 2910.8572
$ws.Range("K126").Value = 8732.571599999999
$ws.Range("M126").Value = -6262.571599999999
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 292.6154
$ws.Range("J12").Value = 310.2
$ws.Range("L12").Value = 930.5999999999999
$ws.Range("N12").Value = -1276.6
$ws.Range("H46").Value = 5405.75
$ws.Range("I46").Value = 6991.1665
$ws.Range("J46").Value = 649.5
$ws.Range("K46").Value = 20973.4995
$ws.Range("L46").Value = 1948.5
$ws.Range("M46").Value = -20882.4995
$ws.Range("N46").Value = -2130.5
$ws.Range("H81").Value = 9648.166999999999
$ws.Range("I81").Value = 1694.5
$ws.Range("K81").Value = 5083.5
$ws.Range("M81").Value = -3960.5
$ws.Range("H84").Value = 9648.166999999999
$ws.Range("I84").Value = 1694.5
$ws.Range("K84").Value = 15250.5
$ws.Range("M84").Value = -9634.5
$ws.Range("H139").Value = 3928.125
$ws.Range("I139").Value = 3775
$ws.Range("K139").Value = 11325
$ws.Range("M139").Value = -6185

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 7000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H131").Value = 23599.2
$ws.Range("I131").Value = 16999
$ws.Range("K131").Value = 16999
$ws.Range("M131").Value = -11959
$ws.Range("H132").Value = 2498.353
$ws.Range("I132").Value = 2246.625
$ws.Range("J132").Value = 2722.111
$ws.Range("K132").Value = 6739.875
$ws.Range("L132").Value = 8166.333
$ws.Range("M132").Value = -4209.875
$ws.Range("N132").Value = -13226.333
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5696.6
$ws.Range("I32").Value = 5327.6665
$ws.Range("J32").Value = 6250
$ws.Range("K32").Value = 5327.6665
$ws.Range("L32").Value = 6250
$ws.Range("M32").Value = -5010.6665
$ws.Range("N32").Value = -6884
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H132").Value = 4185.7856
$ws.Range("I132").Value = 4263.1816
$ws.Range("J132").Value = 3902
$ws.Range("K132").Value = 12789.5448
$ws.Range("L132").Value = 11706
$ws.Range("M132").Value = -10259.5448
$ws.Range("N132").Value = -16766
$ws.Range("H136").Value = 1929.3334
$ws.Range("J136").Value = 1800
$ws.Range("L136").Value = 5400
$ws.Range("N136").Value = -10500
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4929.2
$ws.Range("I132").Value = 5199.1113
$ws.Range("K132").Value = 15597.3339
$ws.Range("M132").Value = -13067.3339
$ws.Range("H136").Value = 7822.923
$ws.Range("I136").Value = 5339.6
$ws.Range("K136").Value = 16018.8
$ws.Range("M136").Value = -13468.8
